$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Split "This sample is compatible with the Windows 10 Creators
#    Update SDK (15063)" into two runs, changing the SDK text to the
#    Fall Creators Update SDK (16299).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Windows 10 Creators Update SDK (15063)", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Windows 10 Fall Creators Update SDK (16299)"
    # Force a run boundary between the untouched text and the freshly
    # typed text (mirrors how Word keeps the edit as its own run).
    $rng.Bold = 1
    $rng.Bold = 0
}

# ------------------------------------------------------------------
# 2. Move the hidden "_GoBack" (last edit position) bookmark from the
#    end of the "October 2017" paragraph to the now-empty paragraph
#    right after the sentence we just edited.
# ------------------------------------------------------------------
$introRange = $d.Content
$introRange.Find.Execute("This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)", `
                          $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$introPara = $introRange.Paragraphs(1)
$nextPara = $introPara.Next()
$d.Bookmarks.Add("_GoBack", $nextPara.Range) | Out-Null

# ------------------------------------------------------------------
# 3. Turn on odd/even headers & footers for the lone section so the
#    blank "default"/"even" header & footer parts are created, and
#    the original table-based header/footer content slides down into
#    the newly created parts (matching the renumbered r:id scheme).
# ------------------------------------------------------------------
$sec = $d.Sections(1)
# wdHeaderFooterEvenPages = 3
$evenHeader = $sec.Headers(3)
$evenHeader.Range.Text = ""
# wdHeaderFooterEvenPages = 3
$evenFooter = $sec.Footers(3)
$evenFooter.Range.Text = ""

Write-Host "Done"
